$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has columns A (segment label) .. E (PercActivationsRescaled)
# for rows 1 (header) .. 20 (data). We need to insert a new column before B that
# holds a numeric "segments" id (0-based row index), push the old label column
# (A) into the new B column, and shift the four numeric data columns
# (old B..E) one slot to the right (new C..F).
# ---------------------------------------------------------------------------

# 1) Shift the four data columns one column to the right (E->F, D->E, C->D,
#    B->C), working right-to-left so we never clobber data we still need.
$ws.Range("E1:E20").Copy($ws.Range("F1:F20"))
$ws.Range("D1:D20").Copy($ws.Range("E1:E20"))
$ws.Range("C1:C20").Copy($ws.Range("D1:D20"))
$ws.Range("B1:B20").Copy($ws.Range("C1:C20"))

# 2) Move the old label column (A) into the new B column. Only the data
#    rows (2-20) need to move - row 1's B1 already holds the old
#    "PercActivations" header cell (moved there in step 1) whose style we
#    want to keep, just with new text. The label cells themselves carry no
#    special formatting (unlike the bold/bordered id column they came
#    from), so the style is reset to Normal after copying the text over.
$ws.Range("B2:B20").Value = $ws.Range("A2:A20").Value()
$ws.Range("B2:B20").Style = "Normal"

# 3) Fix up the header row: B1 becomes "segments", keeping its existing
#    (bold/bordered) style.
$ws.Range("B1").Value = "segments"

# 4) Re-point column A: row 1 has no header label (matches the original,
#    where A1 is blank), rows 2-20 get the numeric 0-based segment id with
#    the bold/bordered style that used to live on the label column.
$ws.Range("A1").ClearContents()

$segmentIds = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18)
for ($i = 0; $i -lt $segmentIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $segmentIds[$i]
}

Write-Output "done"
